# Update the cryptos list with freshly scraped prices / volume changes.
# Price cells (column D) are stored as plain text in the source data (they use
# "." as a thousands separator as well as a decimal point, e.g. "26.497.98"),
# so force a text number format before writing the value to avoid Excel
# re-interpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# Row 2 - Bitcoin
Set-PriceText "D2" "26.497.98"
$ws.Range("E2").Value = "  -1.40%  "

# Row 3 - Ethereum
Set-PriceText "D3" "1.808.45"
$ws.Range("E3").Value = "  -1.13%  "

# Row 4 - TetherUSD
Set-PriceText "D4" "1.005"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5 - now USDC (was BNB)
$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-PriceText "D5" "1.005"
$ws.Range("E5").Value = "  -0.16%  "

# Row 6 - now BNB (was USDC)
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-PriceText "D6" "307.62"
$ws.Range("E6").Value = "  -1.20%  "

# Row 7 - XRP
Set-PriceText "D7" "0.4534"
$ws.Range("E7").Value = "  -1.01%  "

# Row 8 - Cardano
Set-PriceText "D8" "0.3603"
$ws.Range("E8").Value = "  -1.96%  "

# Row 9 - OKB
Set-PriceText "D9" "46.51"
$ws.Range("E9").Value = "  +1.47%  "

# Row 10 - Dogecoin
Set-PriceText "D10" "0.07104"
$ws.Range("E10").Value = "  -0.96%  "

# Row 11 - Polygon
Set-PriceText "D11" "0.8904"
$ws.Range("E11").Value = "  +1.95%  "

# Row 12 - TRON
Set-PriceText "D12" "0.07808"
$ws.Range("E12").Value = "  -0.05%  "

# Row 13 - Solana
Set-PriceText "D13" "19.53"
$ws.Range("E13").Value = "  -0.41%  "

# Row 14 - WrappedEther
Set-PriceText "D14" "1.848.97"
$ws.Range("E14").Value = "  +0.30%  "

# Row 15 - Polkadot
Set-PriceText "D15" "5.300"
$ws.Range("E15").Value = "  -0.59%  "

# Row 16 - Chainlink
Set-PriceText "D16" "6.347"
$ws.Range("E16").Value = "  -0.54%  "

# Row 17 - Litecoin
Set-PriceText "D17" "85.25"
$ws.Range("E17").Value = "  -2.22%  "

# Row 18 - BinanceUSD (only price changes)
Set-PriceText "D18" "1.007"

# Row 19 - ShibaInu
Set-PriceText "D19" "0.000008520"
$ws.Range("E19").Value = "  -2.46%  "

# Row 20 - Dai
Set-PriceText "D20" "1.004"
$ws.Range("E20").Value = "  -0.23%  "

# Row 21 - WrappedBTC
Set-PriceText "D21" "26.526.18"
$ws.Range("E21").Value = "  -1.40%  "

# Row 22 - Avalanche
Set-PriceText "D22" "14.29"
$ws.Range("E22").Value = "  -1.40%  "

# Row 23 - Uniswap
Set-PriceText "D23" "4.981"
$ws.Range("E23").Value = "  -0.12%  "

# Row 24 - WrappedliquidstakedEther2.0
Set-PriceText "D24" "2.053.34"
$ws.Range("E24").Value = "  -1.19%  "

# Row 25 - Cosmos
Set-PriceText "D25" "10.55"
$ws.Range("E25").Value = "  +0.79%  "

# Row 26 - Toncoin
Set-PriceText "D26" "1.975"
$ws.Range("E26").Value = "  -1.35%  "

# Row 27 - Monero
Set-PriceText "D27" "151.14"
$ws.Range("E27").Value = "  -0.10%  "

# Row 28 - EthereumClassic
Set-PriceText "D28" "17.86"
$ws.Range("E28").Value = "  -1.99%  "

# Row 29 - LidoDAOToken
Set-PriceText "D29" "2.061"
$ws.Range("E29").Value = "  +3.96%  "

# Row 30 - BitcoinCash
Set-PriceText "D30" "112.26"
$ws.Range("E30").Value = "  -1.38%  "

# Row 31 - InternetComputer(DFINITY)
Set-PriceText "D31" "4.881"
$ws.Range("E31").Value = "  -1.19%  "

# Row 32 - Stellar
Set-PriceText "D32" "0.08717"
$ws.Range("E32").Value = "  -0.98%  "

# Row 33 - HuobiToken
Set-PriceText "D33" "3.138"
$ws.Range("E33").Value = "  +1.73%  "

# Row 34 - RenderToken
Set-PriceText "D34" "2.854"
$ws.Range("E34").Value = "  +14.15%  "

# Row 35 - Filecoin
Set-PriceText "D35" "4.448"
$ws.Range("E35").Value = "  -0.89%  "

# Row 36 - ImmutableX
Set-PriceText "D36" "0.7245"
$ws.Range("E36").Value = "  -2.92%  "

# Row 37 - ARBITRUM
Set-PriceText "D37" "1.112"
$ws.Range("E37").Value = "  -1.64%  "

# Row 38 - Frax (only volume changes)
$ws.Range("E38").Value = "  -0.27%  "

# Row 39 - TrustWalletToken
Set-PriceText "D39" "1.074"
$ws.Range("E39").Value = "  -0.96%  "

# Row 40 - VeChain
Set-PriceText "D40" "0.01940"
$ws.Range("E40").Value = "  +0.18%  "

# Row 41 - Hedera
Set-PriceText "D41" "0.05107"
$ws.Range("E41").Value = "  -0.37%  "

# Row 42 - MXToken
Set-PriceText "D42" "2.891"
$ws.Range("E42").Value = "  -0.63%  "

# Row 43 - TheSandbox
Set-PriceText "D43" "0.5177"
$ws.Range("E43").Value = "  +4.20%  "

# Row 44 - FraxShare
Set-PriceText "D44" "6.808"
$ws.Range("E44").Value = "  -1.31%  "

# Row 45 - Algorand
Set-PriceText "D45" "0.1515"
$ws.Range("E45").Value = "  -5.04%  "

# Row 46 - Aptos
Set-PriceText "D46" "8.051"
$ws.Range("E46").Value = "  -2.68%  "

# Row 47 - Decentraland
Set-PriceText "D47" "0.4681"
$ws.Range("E47").Value = "  +0.27%  "

# Row 48 - PaxDollar (only volume changes)
$ws.Range("E48").Value = "  -0.20%  "

# Row 49 - EnergySwap
Set-PriceText "D49" "9.970"
$ws.Range("E49").Value = "  -1.50%  "

# Row 50 - Quant
Set-PriceText "D50" "101.34"
$ws.Range("E50").Value = "  -1.57%  "

# Row 51 - NEARProtocol
Set-PriceText "D51" "1.578"
$ws.Range("E51").Value = "  -1.86%  "
